$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.780.60'
$ws.Range('E2').Value = '  +0.38%  '
$ws.Range('D3').Value = '1.915.37'
$ws.Range('E3').Value = '  +1.44%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '241.29'
$ws.Range('E5').Value = '  -2.15%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('E6').Value = '  +0.12%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4910'
$ws.Range('E7').Value = '  -0.36%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2974'
$ws.Range('E8').Value = '  +0.92%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06774'
$ws.Range('E9').Value = '  -0.23%  '
$ws.Range('D10').Value = '1.893.57'
$ws.Range('E10').Value = '  +0.30%  '
$ws.Range('E11').Value = '  -0.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.07372'
$ws.Range('E12').Value = '  +1.79%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.177'
$ws.Range('E13').Value = '  +2.39%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '89.27'
$ws.Range('E14').Value = '  -2.11%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6727'
$ws.Range('E15').Value = '  -0.60%  '
$ws.Range('D16').Value = '30.758.01'
$ws.Range('E16').Value = '  +0.37%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000007986'
$ws.Range('E17').Value = '  +0.24%  '
$ws.Range('E18').Value = '  +2.87%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '1.002'
$ws.Range('E19').Value = '  +0.15%  '
$ws.Range('D20').Value = '2.131.12'
$ws.Range('E20').Value = '  +0.02%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.003'
$ws.Range('E21').Value = '  +0.13%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.284'
$ws.Range('E22').Value = '  +9.63%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '202.99'
$ws.Range('E23').Value = '  +10.82%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.285'
$ws.Range('E24').Value = '  +4.11%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '9.659'
$ws.Range('E25').Value = '  +3.48%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '161.48'
$ws.Range('E26').Value = '  +3.56%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.89'
$ws.Range('E27').Value = '  -0.30%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.970'
$ws.Range('E28').Value = '  +3.63%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.437'
$ws.Range('E29').Value = '  +2.61%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.351'
$ws.Range('E30').Value = '  +1.53%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.09214'
$ws.Range('E31').Value = '  +2.47%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.05568'
$ws.Range('E32').Value = '  +7.29%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '4.082'
$ws.Range('E33').Value = '  +2.38%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.7491'
$ws.Range('E34').Value = '  +1.19%  '
$ws.Range('E35').Value = '  +1.09%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.702'
$ws.Range('E36').Value = '  -1.48%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01864'
$ws.Range('E37').Value = '  +1.62%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9263'
$ws.Range('E39').Value = '  -1.08%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.080'
$ws.Range('E40').Value = '  -3.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.4504'
$ws.Range('E41').Value = '  +2.21%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '72.72'
$ws.Range('E42').Value = '  +25.93%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '107.76'
$ws.Range('E43').Value = '  +2.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.915'
$ws.Range('E44').Value = '  +2.78%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.010'
$ws.Range('E45').Value = '  +0.80%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.1392'
$ws.Range('E46').Value = '  +4.62%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '7.739'
$ws.Range('E47').Value = '  +1.95%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '35.99'
$ws.Range('E48').Value = '  +7.64%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '9.114'
$ws.Range('E49').Value = '  +5.23%  '
$ws.Range('B50').Value = 'Cronos'
$ws.Range('C50').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06006'
$ws.Range('E50').Value = '  +2.80%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4070'
$ws.Range('E51').Value = '  +3.88%  '
